# Update values in the "F" (E column header) and "A" columns as per the
# upstream data correction ("Update Name of Algo" commit).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 12.3073
$ws.Range("E3").Value = 13
$ws.Range("E5").Value = 12.84269999999999

$ws.Range("A9").Value = -20.26489999999997

$ws.Range("E11").Value = 13.31949999999999
$ws.Range("E12").Value = 13.06739999999999

$ws.Range("A13").Value = -21.95820000000002
$ws.Range("A16").Value = -19.93569999999999
$ws.Range("A18").Value = -22.95560000000002
$ws.Range("A20").Value = -22.08910000000003

$ws.Range("E21").Value = 13.10909999999999
